# Pipeline updated for barcode script.
#
# Build-environment marker fix: the docx4j "Modified by ..." banner baked
# into this fixture was generated on a box running Oracle's JDK; the CI
# pipeline now runs the Microsoft Build of OpenJDK, so the banner text
# needs to read "Microsoft Java 21.0.8" instead of "Oracle Java 21.0.8".
#
# Do the substitution the normal Word way: Find & Replace across every
# story in the document (main body, headers/footers, footnotes/endnotes)
# so the fix lands wherever the banner text actually lives, regardless of
# which story happens to hold it.

$d = $word.ActiveDocument

$oldText = "Oracle Java 21.0.8"
$newText = "Microsoft Java 21.0.8"

foreach ($story in $d.StoryRanges) {
    $range = $story
    while ($null -ne $range) {
        $range.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $newText, 2) | Out-Null
        $range = $range.NextStoryRange
    }
}
